# tspi/ciclo-1/task1 -- "Actualizacion del plan general"
#
# 1) Update the activity text in B2 (shared string "Experimento GitHub #1."
#    -> "Ver video tutorial de GitHub.").
# 2) Move the active selection to B2.
# 3) Nudge the column widths of A:F slightly (matches the tiny width
#    recalculation baked into the original commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Ver video tutorial de GitHub."

$ws.Columns.Item(1).ColumnWidth = 11
$ws.Columns.Item(2).ColumnWidth = 39.1666666666667
$ws.Columns.Item(3).ColumnWidth = 1.83333333333333
$ws.Columns.Item(4).ColumnWidth = 11
$ws.Columns.Item(5).ColumnWidth = 1.83333333333333
$ws.Columns.Item(6).ColumnWidth = 19.1666666666667

$ws.Range("B2").Select()
